# Update "Pais" COVID data sheet with refreshed per-country statistics and
# re-sort the data rows (4:218) descending by total cases (column B), matching
# the source data refresh ("Datos actualizados ... a las 20:03").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last updated" timestamp banner -------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 20:03"

# --- 2. Push refreshed statistics into the (pre-sort) country rows ---------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
function Set-Row($r, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

Set-Row 4   1200794 12672 181526 950152 16039 518 69116   # Estados Unidos
Set-Row 11  127659  1614  68166  56032  1384  64  3461    # Turquia
Set-Row 12  102719  1572  42991  52622  8318  81  7106    # Brasil
Set-Row 15  60616   1142  25422  31352  557   160 3842    # Canada
Set-Row 18  44870   2365  12492  30926  0     61  1452    # India
Set-Row 21  29981   76    24500  3697   141   22  1784    # Suiza
Set-Row 27  20941   857   5635   14830  111   19  476     # Pakistan
Set-Row 35  14730   567   2966   11627  1     11  137     # Emiratos Arabes Unidos
Set-Row 57  4783    0     1442   3092   157   3   249     # Argentina
Set-Row 60  4049    129   1173   2847   40    2   29      # Kazajistan
Set-Row 72  2719    550   294    2407   4     0   18      # Ghana
Set-Row 82  1710    124   450    1251   0     2   9       # Guinea
Set-Row 102 742     3     399    337    5     0   6       # Costa Rica
Set-Row 103 751     33    194    549    1     1   8       # Sri Lanka
Set-Row 115 541     14    18     522    2     0   1       # Maldivas
Set-Row 129 325     4     271    31     21    1   23      # Isla de Man
Set-Row 135 236     7     26     200    0     1   10      # Congo
Set-Row 199 15      0     7      7      0     0   1       # Burundi (unchanged numbers)

# --- 3. Re-sort the country table (rows 4-218) by total cases, descending --
$sortRange = $ws.Range("A4:H218")
$sortKey = $ws.Range("B4:B218")
$sortRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)

# --- 4. Fix up the one tie (both 15 total cases) so Burundi precedes ------
#        San Cristobal y Nieves, matching the source refresh ordering.
#        (Excel's stable sort keeps the pre-existing relative order for
#        ties, but the source refresh flipped this particular pair, so we
#        swap the two rows' contents directly, cell by cell.)
function Get-RowVals($r) {
    $vals = @()
    for ($c = 1; $c -le 8; $c++) {
        $vals += $ws.Cells.Item($r, $c).Text
    }
    return $vals
}
function Set-RowVals($r, $vals) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

$rBurundi = 0
$rSanCristobal = 0
for ($r = 4; $r -le 218; $r++) {
    $name = $ws.Cells.Item($r, 1).Text
    if ($name -eq "Burundi") { $rBurundi = $r }
    if ($name -eq "San Cristobal y Nieves") { $rSanCristobal = $r }
}
if ($rBurundi -ne 0 -and $rSanCristobal -ne 0 -and $rBurundi -gt $rSanCristobal) {
    $v1 = Get-RowVals $rSanCristobal
    $v2 = Get-RowVals $rBurundi
    Set-RowVals $rSanCristobal $v2
    Set-RowVals $rBurundi $v1
}
